# Update odds figures in Sheet1 to reflect the latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.25
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 1.13
$ws.Range("G4").Value = 2.5
$ws.Range("I4").Value = 3.4
$ws.Range("J4").Value = 3.5
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("Z4").Value = 2
$ws.Range("AD4").Value = 10
$ws.Range("AF4").Value = 26
$ws.Range("AN4").Value = 6.5
$ws.Range("G5").Value = 1.67
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 5.5
$ws.Range("J5").Value = 2.5
$ws.Range("N5").Value = 4.55
$ws.Range("P5").Value = 2.02
$ws.Range("S5").Value = 3.4
$ws.Range("T5").Value = 1.33
$ws.Range("W5").Value = 5.8
$ws.Range("X5").Value = 1.1
$ws.Range("Y5").Value = 1.73
$ws.Range("Z5").Value = 2.08
$ws.Range("AA5").Value = 2.87
$ws.Range("AB5").Value = 1.36
$ws.Range("AC5").Value = 4
$ws.Range("AN5").Value = 8.5
$ws.Range("AO5").Value = 26
$ws.Range("G7").Value = 1.7
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 5.5
$ws.Range("J7").Value = 2.25
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 5.7
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 5.8
$ws.Range("O7").Value = 1.45
$ws.Range("P7").Value = 2.57
$ws.Range("S7").Value = 2.3
$ws.Range("T7").Value = 1.55
$ws.Range("W7").Value = 4.05
$ws.Range("X7").Value = 1.2
$ws.Range("Y7").Value = 1.47
$ws.Range("Z7").Value = 2.5
$ws.Range("AD7").Value = 6.9
$ws.Range("AF7").Value = 13
$ws.Range("AI7").Value = 5.8
$ws.Range("AJ7").Value = 6.4
$ws.Range("AN7").Value = 11
$ws.Range("AO7").Value = 30
$ws.Range("AP7").Value = 18.5
$ws.Range("AR7").Value = 75
$ws.Range("AS7").Value = 80
$ws.Range("H8").Value = 2.52
$ws.Range("J8").Value = 3.85
$ws.Range("K8").Value = 1.75
$ws.Range("L8").Value = 3.65
$ws.Range("M8").Value = 1.18
$ws.Range("N8").Value = 4.25
$ws.Range("O8").Value = 1.75
$ws.Range("P8").Value = 1.98
$ws.Range("S8").Value = 3.15
$ws.Range("T8").Value = 1.31
$ws.Range("W8").Value = 5.9
$ws.Range("X8").Value = 1.1
$ws.Range("Y8").Value = 1.7
$ws.Range("Z8").Value = 2.02
$ws.Range("AA8").Value = 2.37
$ws.Range("AB8").Value = 1.52
$ws.Range("AC8").Value = 5.9
$ws.Range("AD8").Value = 13.5
$ws.Range("AH8").Value = 70
$ws.Range("AI8").Value = 4.25
$ws.Range("AJ8").Value = 5.4
$ws.Range("AK8").Value = 21
$ws.Range("AL8").Value = 175
$ws.Range("AP8").Value = 11.75
$ws.Range("AS8").Value = 65
$ws.Range("G9").Value = 1.39
$ws.Range("I9").Value = 8.5
$ws.Range("K9").Value = 2.25
$ws.Range("L9").Value = 7.1
$ws.Range("N9").Value = 8
$ws.Range("O9").Value = 1.23
$ws.Range("P9").Value = 3.7
$ws.Range("S9").Value = 1.72
$ws.Range("T9").Value = 2.02
$ws.Range("W9").Value = 2.67
$ws.Range("X9").Value = 1.42
$ws.Range("Y9").Value = 1.37
$ws.Range("Z9").Value = 2.85
$ws.Range("AA9").Value = 1.93
$ws.Range("AB9").Value = 1.78
$ws.Range("AC9").Value = 6.6
$ws.Range("AD9").Value = 6.5
$ws.Range("AE9").Value = 8
$ws.Range("AF9").Value = 9
$ws.Range("AG9").Value = 11.25
$ws.Range("AH9").Value = 26
$ws.Range("AI9").Value = 8
$ws.Range("AK9").Value = 18
$ws.Range("AL9").Value = 80
$ws.Range("AM9").Value = 600
$ws.Range("AN9").Value = 23
$ws.Range("AP9").Value = 24
$ws.Range("AS9").Value = 70
$ws.Range("G10").Value = 1.82
$ws.Range("I10").Value = 3.4
$ws.Range("O10").Value = 1.14
$ws.Range("P10").Value = 5.5
$ws.Range("S10").Value = 1.53
$ws.Range("T10").Value = 2.4
$ws.Range("U10").Value = 1.8
$ws.Range("V10").Value = 2
$ws.Range("W10").Value = 2.2
$ws.Range("X10").Value = 1.62
$ws.Range("AA10").Value = 1.53
$ws.Range("AB10").Value = 2.38
$ws.Range("AD10").Value = 11
$ws.Range("AF10").Value = 17
$ws.Range("AG10").Value = 13
$ws.Range("AR10").Value = 26
$ws.Range("S12").Value = 1.5
$ws.Range("T12").Value = 2.5
$ws.Range("AM12").Value = 1000
$ws.Range("S13").Value = 1.57
$ws.Range("T13").Value = 2.35
$ws.Range("U13").Value = 1.95
$ws.Range("V13").Value = 1.85
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = 4.1
$ws.Range("I14").Value = 1.55
$ws.Range("K14").Value = 2.25
$ws.Range("AS14").Value = 26
$ws.Range("M15").Value = 1.02
$ws.Range("N15").Value = 19
$ws.Range("O15").Value = 1.14
$ws.Range("P15").Value = 5.5
$ws.Range("S15").Value = 1.5
$ws.Range("T15").Value = 2.5
$ws.Range("U15").Value = 1.83
$ws.Range("V15").Value = 2.03
$ws.Range("G16").Value = 3.6
$ws.Range("H16").Value = 3.5
$ws.Range("I16").Value = 2.05
$ws.Range("L16").Value = 2.63
$ws.Range("M16").Value = 1.04
$ws.Range("N16").Value = 13
$ws.Range("S16").Value = 1.83
$ws.Range("T16").Value = 2.03
$ws.Range("AC16").Value = 11
$ws.Range("AE16").Value = 12
$ws.Range("AG16").Value = 26
$ws.Range("AI16").Value = 11
$ws.Range("AJ16").Value = 6.5
$ws.Range("AM16").Value = 201
$ws.Range("AQ16").Value = 19
$ws.Range("G17").Value = 5.3
$ws.Range("H17").Value = 4.05
$ws.Range("I17").Value = 1.55
$ws.Range("J17").Value = 5.1
$ws.Range("K17").Value = 2.32
$ws.Range("L17").Value = 2.07
$ws.Range("N17").Value = 8.25
$ws.Range("O17").Value = 1.22
$ws.Range("P17").Value = 3.85
$ws.Range("T17").Value = 2.1
$ws.Range("W17").Value = 2.6
$ws.Range("X17").Value = 1.44
$ws.Range("Y17").Value = 1.33
$ws.Range("Z17").Value = 3.05
$ws.Range("AA17").Value = 1.75
$ws.Range("AB17").Value = 1.98
$ws.Range("AC17").Value = 16.5
$ws.Range("AD17").Value = 32
$ws.Range("AE17").Value = 16.5
$ws.Range("AF17").Value = 100
$ws.Range("AG17").Value = 50
$ws.Range("AH17").Value = 45
$ws.Range("AI17").Value = 8.25
$ws.Range("AJ17").Value = 7.9
$ws.Range("AK17").Value = 15.5
$ws.Range("AL17").Value = 65
$ws.Range("AM17").Value = 450
$ws.Range("AO17").Value = 7.8
$ws.Range("AQ17").Value = 11.25
$ws.Range("AR17").Value = 11.75
$ws.Range("I18").Value = 3.8
$ws.Range("L18").Value = 4.2
$ws.Range("S18").Value = 1.83
$ws.Range("W18").Value = 2.92
$ws.Range("X18").Value = 1.3
$ws.Range("AA18").Value = 1.78
$ws.Range("AB18").Value = 1.83
$ws.Range("AC18").Value = 7.1
$ws.Range("AE18").Value = 8.5
$ws.Range("AF18").Value = 14.5
$ws.Range("AI18").Value = 10.25
$ws.Range("AJ18").Value = 7
$ws.Range("AL18").Value = 75
$ws.Range("AO18").Value = 20
$ws.Range("AS18").Value = 45